# Edit the "Requisitos" list in the LOM3210 document.
#
# The source list (one course requirement per run, each run holding
# "<w:t>text</w:t><w:br/>") is reordered, three items are dropped, and
# one item is added, with an accent fix on "Álgebra". We rebuild the
# paragraph by matching the *old* 25-item sequence positionally against
# the *new* 23-item sequence (the first 23 slots are updated in place,
# preserving run/formatting identity; the trailing 2 slots are removed).

$oldItems = @(
    "LOM3236 -  Processos de Fabricação  (Requisito)",
    "LOB1053 -  Física III  (Requisito)",
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)",
    "LOB1036 -  Geometria Analítica  (Requisito)",
    "LOB1037 -  Àlgebra Linear  (Requisito)",
    "LOB1042 -  Física Experimental IV  (Requisito)",
    "LOB1039 -  Física Experimental III  (Requisito)",
    "LOB1018 -  Física I  (Requisito)",
    "LOB1021 -  Física IV  (Requisito)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
    "LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)",
    "LOB1012 -  Estatística  (Requisito)",
    "LOB1038 -  Física Experimental I  (Requisito)",
    "LOB1052 -  Cálculo III  (Requisito)",
    "LOB1008 -  Ciência, Tecnologia e Sociedade  (Requisito)",
    "LOB1041 -  Física Experimental II  (Requisito)",
    "LOM3240 -  Química Inorgânica Fundamental e Aplicada  (Requisito)",
    "LOQ4095 -  Química Geral Experimental  (Requisito)",
    "LOM3218 -  Introdução à Engenharia Física  (Requisito)",
    "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)",
    "LOB1004 -  Cálculo II  (Requisito)",
    "LOB1003 -  Cálculo I  (Requisito)",
    "LOB1006 -  Cálculo IV  (Requisito)",
    "LOM3260 -  Computação Científica em Python  (Requisito)",
    "LOB1019 -  Física II  (Requisito)"
)

$newItems = @(
    "LOQ4095 -  Química Geral Experimental  (Requisito)",
    "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)",
    "LOB1006 -  Cálculo IV  (Requisito)",
    "LOB1053 -  Física III  (Requisito)",
    "LOB1003 -  Cálculo I  (Requisito)",
    "LOB1012 -  Estatística  (Requisito)",
    "LOB1036 -  Geometria Analítica  (Requisito)",
    "LOB1037 -  Álgebra Linear  (Requisito)",
    "LOB1038 -  Física Experimental I  (Requisito)",
    "LOB1039 -  Física Experimental III  (Requisito)",
    "LOB1041 -  Física Experimental II  (Requisito)",
    "LOB1042 -  Física Experimental IV  (Requisito)",
    "LOB1052 -  Cálculo III  (Requisito)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
    "LOM3037 -  Química Inorgânica  (Requisito)",
    "LOM3260 -  Computação Científica em Python  (Requisito)",
    "LOB1004 -  Cálculo II  (Requisito)",
    "LOB1018 -  Física I  (Requisito)",
    "LOB1019 -  Física II  (Requisito)",
    "LOB1021 -  Física IV  (Requisito)",
    "LOM3218 -  Introdução à Engenharia Física  (Requisito)",
    "LOM3204 -  Desenho Técnico e Projeto Assistido por Computador  (Requisito)",
    "LOM3236 -  Processos de Fabricação  (Requisito)"
)

$d = $word.ActiveDocument

# Locate the requirements list paragraph: the ListBullet paragraph that
# immediately follows the paragraph containing "Requisitos".
$reqParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt.TrimEnd([char]13, [char]7) -eq "Requisitos") {
        $reqParaIndex = $i + 1
        break
    }
}
if ($reqParaIndex -eq -1) { throw "Could not find the Requisitos heading paragraph" }

$p = $d.Paragraphs.Item($reqParaIndex)
$paraStart = $p.Range.Start

# Compute the start/end offsets (relative to $paraStart) of each of the
# existing items by splitting on the manual line break character used
# between runs (Chr(11), i.e. <w:br/>).
$fullText = $p.Range.Text
$bodyText = $fullText.Substring(0, $fullText.Length - 1)  # drop trailing paragraph mark
$parts = $bodyText -split [char]11

if ($parts.Count -ne ($oldItems.Count + 1)) {
    throw ("Unexpected item count in Requisitos paragraph: " + $parts.Count)
}

$starts = @()
$ends = @()
$pos = 0
for ($i = 0; $i -lt $oldItems.Count; $i++) {
    $len = $parts[$i].Length
    $starts += $pos
    $ends += ($pos + $len)
    $pos = $pos + $len + 1
}
$paraEndOffset = $pos

# Verify the existing runs match the expected original text before editing.
for ($i = 0; $i -lt $oldItems.Count; $i++) {
    $r = $d.Range($paraStart + $starts[$i], $paraStart + $ends[$i])
    if ($r.Text -ne $oldItems[$i]) {
        throw ("Mismatch at item " + $i + ": expected '" + $oldItems[$i] + "' got '" + $r.Text + "'")
    }
}

# Drop the trailing items that have no counterpart in the new list
# (remove from the end first so earlier offsets stay valid).
if ($newItems.Count -lt $oldItems.Count) {
    $cutStart = $starts[$newItems.Count]
    $rCut = $d.Range($paraStart + $cutStart, $paraStart + $paraEndOffset)
    $rCut.Text = ""
}

# Update the remaining runs' text in place, last-to-first so offsets
# computed above stay valid as the document shrinks/grows.
for ($i = $newItems.Count - 1; $i -ge 0; $i--) {
    $r = $d.Range($paraStart + $starts[$i], $paraStart + $ends[$i])
    $r.Text = $newItems[$i]
}

Write-Output "Requisitos list updated successfully"
